$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 169, shifting rows 169:201 down to 170:202
$ws.Rows.Item(169).Insert()

# Fill in the new row 169 with the new data record
$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 44476
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100112044
$ws.Cells.Item(169, 7).Value = "Perejil"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 40
$ws.Cells.Item(169, 11).Value = 4000
$ws.Cells.Item(169, 12).Value = 4000
$ws.Cells.Item(169, 13).Value = 4000
$ws.Cells.Item(169, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(169, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(169, 16).Value = 1333
$ws.Cells.Item(169, 17).Value = 3
$ws.Cells.Item(169, 18).Value = "Hortaliza"
